# Add two new columns (Role, City) to the existing data table and
# update the age for the last row (伊藤美咲) from 27 to 22.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C: Role (職種)
$ws.Range("C1").Value = "エンジニア"
$ws.Range("C2").Value = "デザイナー"
$ws.Range("C3").Value = "営業"
$ws.Range("C4").Value = "マネージャー"
$ws.Range("C5").Value = "デザイナー"

# Column D: City (都市)
$ws.Range("D1").Value = "東京"
$ws.Range("D2").Value = "大阪"
$ws.Range("D3").Value = "福岡"
$ws.Range("D4").Value = "名古屋"
$ws.Range("D5").Value = "横浜"

# Update age for row 5 (伊藤美咲) from 27 to 22
$ws.Range("B5").Value = 22
